# OPTIMIZACION DE DISEÑOS Y PROCESOS NO.13
# Corrige el estado de Carlos Perez (INACTIVO -> ACTIVO) y actualiza el
# correo electronico de Jose Ramirez a su cuenta de gmail.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "ACTIVO"
$ws.Range("E7").Value = "jose.ramirez@gmail.com"

$ws.Range("L7").Select()
